$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was collected for "Feria Lagunitas de Puerto
# Montt - Membrillo". It belongs at the top of the date-ordered block of
# data rows (row 91), so push the existing rows 91:118 down by one and
# fill the freed row with the new record.
$ws.Rows.Item(91).Insert()

$row = 91
$ws.Cells.Item($row,1).Value  = 4
$ws.Cells.Item($row,2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item($row,3).Value  = "Los Lagos"
$ws.Cells.Item($row,4).Value  = 44985
$ws.Cells.Item($row,5).Value  = 10
$ws.Cells.Item($row,6).Value  = "Fruta"
$ws.Cells.Item($row,7).Value  = 100104
$ws.Cells.Item($row,8).Value  = "Frutos de pepita"
$ws.Cells.Item($row,9).Value  = 100104003
$ws.Cells.Item($row,10).Value = "Membrillo"
$ws.Cells.Item($row,11).Value = "Champion"
$ws.Cells.Item($row,12).Value = "Primera"
$ws.Cells.Item($row,13).Value = 400
$ws.Cells.Item($row,14).Value = 17000
$ws.Cells.Item($row,15).Value = 18000
$ws.Cells.Item($row,16).Value = 17500
$ws.Cells.Item($row,17).Value = "`$/caja 18 kilos empedrada"
$ws.Cells.Item($row,18).Value = "Región de O'Higgins"
$ws.Cells.Item($row,19).Value = 972
$ws.Cells.Item($row,20).Value = 18
